$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "309.83"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "-3.44%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "50.70"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "3.58%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.172"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-1.47%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.07782"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "4.497"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "-2.15%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.354"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "11.97%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.567"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-4.45%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.1210"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-6.25%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1976"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "1.91%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.04740"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "2.60%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.09375"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-0.34%"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-0.68%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.001259"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-5.35%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.005785"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-0.87%"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "2,012.59%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.329"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-0.32%"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.434"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "2.04%"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.992"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "-1.12%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.1371"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "-0.79%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04166"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-0.25%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001271"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-2.78%"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-7.01%"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0001350"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02604"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "-4.06%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.06008"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "5.57%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01100"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "74.22%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007842"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-0.68%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1426"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "-1.24%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.008394"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "8.86%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.007661"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-5.44%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.3384"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "5.96%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00007338"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "6.29%"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000750"
$ws.Range("B48").NumberFormat = "@"
$ws.Range("B48").Value = "BOLO"
$ws.Range("C48").NumberFormat = "@"
$ws.Range("C48").Value = "https://coinranking.com/coin/ogrGe0dEab+bolo-bolo"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.05317"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "-22.86%"
$ws.Range("B49").NumberFormat = "@"
$ws.Range("B49").Value = "CoinbaseStockToken"
$ws.Range("C49").NumberFormat = "@"
$ws.Range("C49").Value = "https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.002619"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "-34.64%"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002100"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0002000"
